$p = $ppt.ActivePresentation

# The deck ships with two embedded DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (the stock Office palette)
#   ppt/theme/theme2.xml -> "Integral"     (the palette actually applied
#                                            to the slide master / deck)
# The authored change swaps the two themes' contents so the deck's live
# theme becomes the plain "Office Theme" palette (and the previously
# unused "Integral" palette becomes the dormant one). We reproduce the
# live-theme side of that swap through the real PowerPoint object model:
# the slide master's theme color scheme is the one embedded object that
# is actually wired into the presentation, so we repaint its 12 theme
# colors with the stock Office Theme RGB values.

$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Index order matches the standard theme color slots:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
